# Auto-generated script applying numeric corrections to Leve profit tables
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 333.33334
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()

$ws.Range("H88").Value = 4555.154
$ws.Range("I88").Value = 2699.6
$ws.Range("K88").Value = 2699.6
$ws.Range("M88").Value = -2293.6

$ws.Range("H91").Value = 4555.154
$ws.Range("I91").Value = 2699.6
$ws.Range("K91").Value = 2699.6
$ws.Range("M91").Value = -1295.6

$ws.Range("H112").Value = 2433.238
$ws.Range("J112").Value = 2504.9
$ws.Range("L112").Value = 7514.700000000001
$ws.Range("N112").Value = -9730.700000000001

$ws.Range("H115").Value = 1737
$ws.Range("I115").Value = 892.5
$ws.Range("J115").Value = 2300
$ws.Range("K115").Value = 2677.5
$ws.Range("L115").Value = 6900
$ws.Range("M115").Value = -1110.5
$ws.Range("N115").Value = -10034

$ws.Range("H135").Value = 362.2857

$ws.Range("H137").Value = 1671451.1
$ws.Range("I137").Value = 4048.4856
$ws.Range("J137").Value = 4005815
$ws.Range("K137").Value = 12145.4568
$ws.Range("L137").Value = 12017445
$ws.Range("M137").Value = -9595.4568
$ws.Range("N137").Value = -12022545

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3304.111
$ws.Range("I32").Value = 1411.2354
$ws.Range("J32").Value = 11348.833
$ws.Range("K32").Value = 1411.2354
$ws.Range("L32").Value = 11348.833
$ws.Range("M32").Value = -1124.2354
$ws.Range("N32").Value = -11922.833

$ws.Range("H36").Value = 9341.333000000001
$ws.Range("I36").Value = 9341.333000000001
$ws.Range("K36").Value = 9341.333000000001
$ws.Range("M36").Value = -8995.333000000001

$ws.Range("H45").Value = 86794
$ws.Range("I45").Value = 86794
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 86794
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -86417
$ws.Range("N45").ClearContents()

$ws.Range("H122").Value = 2134.0625
$ws.Range("I122").Value = 1917.2084
$ws.Range("J122").Value = 2784.625
$ws.Range("K122").Value = 5751.6252
$ws.Range("L122").Value = 8353.875
$ws.Range("M122").Value = -3301.6252
$ws.Range("N122").Value = -13253.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7333
$ws.Range("J86").Value = 15999
$ws.Range("L86").Value = 15999
$ws.Range("N86").Value = -18245

$ws.Range("H89").Value = 7333
$ws.Range("J89").Value = 15999
$ws.Range("L89").Value = 79995
$ws.Range("N89").Value = -91227

$ws.Range("H105").Value = 7287.8887
$ws.Range("I105").Value = 9069.625
$ws.Range("K105").Value = 9069.625
$ws.Range("M105").Value = -7322.625

$ws.Range("H111").Value = 25534.5
$ws.Range("J111").Value = 25534.5
$ws.Range("L111").Value = 25534.5
$ws.Range("N111").Value = -33714.5

$ws.Range("H134").Value = 30002416
$ws.Range("I134").Value = 2351.926
$ws.Range("K134").Value = 7055.778
$ws.Range("M134").Value = -4520.778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 754113.4
$ws.Range("J3").Value = 4500
$ws.Range("L3").Value = 4500
$ws.Range("N3").Value = -4726

$ws.Range("H16").Value = 7529732.5
$ws.Range("I16").Value = 17863192
$ws.Range("K16").Value = 17863192
$ws.Range("M16").Value = -17862905

$ws.Range("H31").Value = 2060.06
$ws.Range("I31").Value = 2159.1428
$ws.Range("J31").Value = 2021.5278
$ws.Range("K31").Value = 2159.1428
$ws.Range("L31").Value = 2021.5278
$ws.Range("M31").Value = -1864.1428
$ws.Range("N31").Value = -2611.5278

$ws.Range("H34").Value = 2060.06
$ws.Range("I34").Value = 2159.1428
$ws.Range("J34").Value = 2021.5278
$ws.Range("K34").Value = 2159.1428
$ws.Range("L34").Value = 2021.5278
$ws.Range("M34").Value = -1957.1428
$ws.Range("N34").Value = -2425.5278

$ws.Range("H113").Value = 7529732.5
$ws.Range("I113").Value = 17863192
$ws.Range("K113").Value = 17863192
$ws.Range("M113").Value = -17861022

$ws.Range("H132").Value = 11907806
$ws.Range("I132").Value = 2085.3333
$ws.Range("J132").Value = 25645176
$ws.Range("K132").Value = 6255.999899999999
$ws.Range("L132").Value = 76935528
$ws.Range("M132").Value = -3725.999899999999
$ws.Range("N132").Value = -76940588

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 2127.5715
$ws.Range("J32").Value = 2315.6667
$ws.Range("L32").Value = 6947.000100000001
$ws.Range("N32").Value = -7513.000100000001

$ws.Range("H68").Value = 1068.4138
$ws.Range("J68").Value = 1040.1482
$ws.Range("L68").Value = 3120.4446
$ws.Range("N68").Value = -4742.444600000001

$ws.Range("H71").Value = 1068.4138
$ws.Range("J71").Value = 1040.1482
$ws.Range("L71").Value = 9361.3338
$ws.Range("N71").Value = -17473.3338

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 2504
$ws.Range("J6").Value = 2504
$ws.Range("L6").Value = 2504
$ws.Range("N6").Value = -2730

$ws.Range("H10").Value = 1743422.5
$ws.Range("I10").Value = 3000
$ws.Range("K10").Value = 3000
$ws.Range("M10").Value = -2831

$ws.Range("H13").Value = 1950
$ws.Range("I13").Value = 100
$ws.Range("J13").Value = 2875
$ws.Range("K13").Value = 100
$ws.Range("L13").Value = 2875
$ws.Range("M13").Value = 39
$ws.Range("N13").Value = -3153

$ws.Range("H16").Value = 2504
$ws.Range("J16").Value = 2504
$ws.Range("L16").Value = 2504
$ws.Range("N16").Value = -3004

$ws.Range("H19").Value = 500200
$ws.Range("I19").Value = 400
$ws.Range("J19").Value = 1000000
$ws.Range("K19").Value = 400
$ws.Range("L19").Value = 1000000
$ws.Range("M19").Value = -112
$ws.Range("N19").Value = -1000576

$ws.Range("H107").Value = 64590
$ws.Range("I107").Value = 334220
$ws.Range("J107").Value = 2367.6924
$ws.Range("K107").Value = 334220
$ws.Range("L107").Value = 2367.6924
$ws.Range("M107").Value = -332300
$ws.Range("N107").Value = -6207.6924

$ws.Range("H126").Value = 2564.6667
$ws.Range("I126").Value = 1956.0834
$ws.Range("J126").Value = 4999
$ws.Range("K126").Value = 5868.2502
$ws.Range("L126").Value = 14997
$ws.Range("M126").Value = -3398.2502
$ws.Range("N126").Value = -19937

$ws.Range("H132").Value = 4523007
$ws.Range("I132").Value = 4649.08
$ws.Range("K132").Value = 13947.24
$ws.Range("M132").Value = -11417.24

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8369.591
$ws.Range("I7").Value = 4759.923
$ws.Range("K7").Value = 4759.923
$ws.Range("M7").Value = -4647.923

$ws.Range("H22").Value = 5907.636
$ws.Range("J22").Value = 7748.125
$ws.Range("L22").Value = 7748.125
$ws.Range("N22").Value = -8338.125

$ws.Range("H27").Value = 5907.636
$ws.Range("J27").Value = 7748.125
$ws.Range("L27").Value = 7748.125
$ws.Range("N27").Value = -7962.125

$ws.Range("H68").Value = 2736.75
$ws.Range("I68").Value = 2875
$ws.Range("J68").Value = 2598.5
$ws.Range("K68").Value = 2875
$ws.Range("L68").Value = 2598.5
$ws.Range("M68").Value = -2126
$ws.Range("N68").Value = -4096.5

$ws.Range("H71").Value = 2736.75
$ws.Range("I71").Value = 2875
$ws.Range("J71").Value = 2598.5
$ws.Range("K71").Value = 14375
$ws.Range("L71").Value = 12992.5
$ws.Range("M71").Value = -10631
$ws.Range("N71").Value = -20480.5

$ws.Range("H82").Value = 1659.6428
$ws.Range("I82").Value = 1543.5
$ws.Range("K82").Value = 1543.5
$ws.Range("M82").Value = -1182.5

$ws.Range("H85").Value = 1659.6428
$ws.Range("I85").Value = 1543.5
$ws.Range("K85").Value = 1543.5
$ws.Range("M85").Value = -295.5

$ws.Range("H122").Value = 3161.8333
$ws.Range("I122").Value = 2995.6365
$ws.Range("K122").Value = 8986.9095
$ws.Range("M122").Value = -6536.9095

$ws.Range("H126").Value = 8369.591
$ws.Range("I126").Value = 4759.923
$ws.Range("K126").Value = 14279.769
$ws.Range("M126").Value = -11809.769

$ws.Range("H132").Value = 3894.122
$ws.Range("I132").Value = 3595.0417
$ws.Range("J132").Value = 4316.353
$ws.Range("K132").Value = 10785.1251
$ws.Range("L132").Value = 12949.059
$ws.Range("M132").Value = -8255.125100000001
$ws.Range("N132").Value = -18009.059

$ws.Range("H136").Value = 2644.12
$ws.Range("I136").Value = 2268.111
$ws.Range("K136").Value = 6804.333
$ws.Range("M136").Value = -4254.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 64443.8
$ws.Range("J94").Value = 64443.8
$ws.Range("L94").Value = 64443.8
$ws.Range("N94").Value = -66245.8
